$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.245.65'
$ws.Cells.Item(2, 5).Value = '  -4.38%  '

$ws.Cells.Item(3, 4).Value = '3.089.63'
$ws.Cells.Item(3, 5).Value = '  -4.26%  '

$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '541.26'
$ws.Cells.Item(5, 5).Value = '  -6.30%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '133.24'
$ws.Cells.Item(6, 5).Value = '  -12.13%  '

$ws.Cells.Item(7, 5).Value = '  +0.05%  '

$ws.Cells.Item(8, 4).Value = '3.083.87'
$ws.Cells.Item(8, 5).Value = '  -4.29%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.489'
$ws.Cells.Item(9, 5).Value = '  -4.53%  '

$ws.Cells.Item(10, 5).Value = '  -4.80%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '6.16'
$ws.Cells.Item(11, 5).Value = '  -12.90%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.460'
$ws.Cells.Item(12, 5).Value = '  -5.67%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.0000226'
$ws.Cells.Item(13, 5).Value = '  -3.14%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '34.51'
$ws.Cells.Item(14, 5).Value = '  -10.27%  '

$ws.Cells.Item(15, 4).Value = '3.540.57'
$ws.Cells.Item(15, 5).Value = '  -5.67%  '

$ws.Cells.Item(16, 4).Value = '63.109.85'
$ws.Cells.Item(16, 5).Value = '  -4.63%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.111'
$ws.Cells.Item(17, 5).Value = '  -3.37%  '

$ws.Cells.Item(18, 4).Value = '3.082.79'
$ws.Cells.Item(18, 5).Value = '  -4.61%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '6.61'
$ws.Cells.Item(19, 5).Value = '  -7.43%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '484.69'
$ws.Cells.Item(20, 5).Value = '  -10.18%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '13.32'
$ws.Cells.Item(21, 5).Value = '  -8.90%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.704'
$ws.Cells.Item(22, 5).Value = '  -5.46%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '7.20'
$ws.Cells.Item(23, 5).Value = '  -6.87%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '78.61'
$ws.Cells.Item(24, 5).Value = '  -3.19%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '12.06'
$ws.Cells.Item(25, 5).Value = '  -10.86%  '

$ws.Cells.Item(26, 5).Value = '  -0.02%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '2.69'
$ws.Cells.Item(27, 5).Value = '  -8.91%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '8.17'
$ws.Cells.Item(28, 5).Value = '  -13.39%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.997'
$ws.Cells.Item(29, 5).Value = '  -0.35%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '26.00'
$ws.Cells.Item(30, 5).Value = '  -5.86%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.90'
$ws.Cells.Item(31, 5).Value = '  -16.37%  '

$ws.Cells.Item(32, 5).Value = '  -6.61%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '58.89'
$ws.Cells.Item(33, 5).Value = '  +7.60%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '2.42'
$ws.Cells.Item(34, 5).Value = '  -12.16%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '5.97'
$ws.Cells.Item(35, 5).Value = '  -6.29%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '5.21'
$ws.Cells.Item(36, 5).Value = '  -7.32%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '464.73'
$ws.Cells.Item(37, 5).Value = '  -17.33%  '

$ws.Cells.Item(38, 4).Value = '3.120.68'
$ws.Cells.Item(38, 5).Value = '  -2.99%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.0390'
$ws.Cells.Item(39, 5).Value = '  -14.56%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.0789'
$ws.Cells.Item(40, 5).Value = '  -8.17%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.115'
$ws.Cells.Item(41, 5).Value = '  -11.91%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '8.06'
$ws.Cells.Item(42, 5).Value = '  -6.34%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '2.53'
$ws.Cells.Item(43, 5).Value = '  -12.74%  '

$ws.Cells.Item(44, 2).Value = 'TheGraph'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.250'
$ws.Cells.Item(44, 5).Value = '  -11.76%  '

$ws.Cells.Item(45, 2).Value = 'USDe'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.999'
$ws.Cells.Item(45, 5).Value = '  +0.01%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '2.03'
$ws.Cells.Item(46, 5).Value = '  -13.67%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '24.50'
$ws.Cells.Item(47, 5).Value = '  -7.41%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '117.77'
$ws.Cells.Item(48, 5).Value = '  -5.66%  '

$ws.Cells.Item(49, 5).Value = '  -4.84%  '

$ws.Cells.Item(50, 4).Value = '0.0₃0512'
$ws.Cells.Item(50, 5).Value = '  -7.65%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.99'
$ws.Cells.Item(51, 5).Value = '  -10.16%  '
